$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.443.55'
$ws.Range("E2").Value = '  -3.64%  '
$ws.Range("D3").Value = '1.993.75'
$ws.Range("E3").Value = '  -6.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.45'
$ws.Range("E5").Value = '  -5.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5009'
$ws.Range("E7").Value = '  -4.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4220'
$ws.Range("E8").Value = '  -5.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.80'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08902'
$ws.Range("E10").Value = '  -5.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.122'
$ws.Range("E11").Value = '  -5.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.21'
$ws.Range("E12").Value = '  -8.34%  '
$ws.Range("D13").Value = '2.017.78'
$ws.Range("E13").Value = '  -3.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.062'
$ws.Range("E14").Value = '  -7.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.511'
$ws.Range("E15").Value = '  -7.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.01'
$ws.Range("E16").Value = '  -6.42%  '
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001109'
$ws.Range("E18").Value = '  -5.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06625'
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.66'
$ws.Range("E20").Value = '  -9.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.008'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.968'
$ws.Range("E22").Value = '  -5.94%  '
$ws.Range("D23").Value = '29.483.39'
$ws.Range("E23").Value = '  -3.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.88'
$ws.Range("E24").Value = '  -7.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.247'
$ws.Range("E25").Value = '  -3.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.90'
$ws.Range("E26").Value = '  -2.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.59'
$ws.Range("E27").Value = '  -7.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.474'
$ws.Range("E28").Value = '  -7.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.336'
$ws.Range("E29").Value = '  -8.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.91'
$ws.Range("E30").Value = '  -5.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.046'
$ws.Range("E31").Value = '  -10.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09929'
$ws.Range("E32").Value = '  -6.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.567'
$ws.Range("E33").Value = '  -12.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.841'
$ws.Range("E34").Value = '  -7.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.789'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.567'
$ws.Range("E36").Value = '  -10.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02457'
$ws.Range("E37").Value = '  -7.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06349'
$ws.Range("E38").Value = '  -7.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.287'
$ws.Range("E39").Value = '  -3.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6504'
$ws.Range("E40").Value = '  -8.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.71'
$ws.Range("E41").Value = '  -7.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2066'
$ws.Range("E42").Value = '  -8.63%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6335'
$ws.Range("E44").Value = '  -8.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.39'
$ws.Range("E45").Value = '  -8.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.200'
$ws.Range("E46").Value = '  -8.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.293'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.523'
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000327'
$ws.Range("E49").Value = '  -5.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06998'
$ws.Range("E50").Value = '  -3.41%  '
$ws.Range("E51").Value = '  -5.33%  '
